$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(45863.01041666666,45863.02083333334,45863.03125,45863.04166666666,45863.05208333334,45863.0625,45863.07291666666,45863.08333333334,45863.09375,45863.10416666666,45863.11458333334,45863.125,45863.13541666666,45863.14583333334,45863.15625,45863.16666666666,45863.17708333334,45863.1875,45863.19791666666,45863.20833333334,45863.21875,45863.22916666666,45863.23958333334,45863.25,45863.26041666666,45863.27083333334,45863.28125,45863.29166666666,45863.30208333334,45863.3125,45863.32291666666,45863.33333333334,45863.34375,45863.35416666666,45863.36458333334,45863.375,45863.38541666666,45863.39583333334,45863.40625,45863.41666666666,45863.42708333334,45863.4375,45863.44791666666,45863.45833333334,45863.46875,45863.47916666666,45863.48958333334,45863.5,45863.51041666666,45863.52083333334,45863.53125,45863.54166666666,45863.55208333334,45863.5625,45863.57291666666,45863.58333333334,45863.59375,45863.60416666666,45863.61458333334,45863.625,45863.63541666666,45863.64583333334,45863.65625,45863.66666666666,45863.67708333334,45863.6875,45863.69791666666,45863.70833333334,45863.71875,45863.72916666666,45863.73958333334,45863.75,45863.76041666666,45863.77083333334,45863.78125,45863.79166666666,45863.80208333334,45863.8125,45863.82291666666,45863.83333333334,45863.84375,45863.85416666666,45863.86458333334,45863.875,45863.88541666666,45863.89583333334,45863.90625,45863.91666666666,45863.92708333334,45863.9375,45863.94791666666,45863.95833333334,45863.96875,45863.97916666666,45863.98958333334,45864)
$values = @(353,357,356,354,326,323,320,319,275,271,270,269,222,221,219,218,200,199,198,195,178,177,177,177,164,163,163,162,135,135,136,137,113,114,116,118,158,160,162,163,216,217,219,221,270,271,271,272,342,342,343,343,415,416,417,418,518,520,521,523,621,623,625,627,703,705,707,708,746,746,746,746,752,751,750,749,777,781,784,788,885,887,889,892,911,914,917,920,903,904,905,905,0,0,0,0)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

"Updated $($dates.Length) rows"
